$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Append 5 new time-log rows (90-94), matching the formatting used by
#     the existing rows right above them (row 89: text in B, plain number
#     in C, date-formatted number in D). ---

# Copy the date number-format (column D) down to the new rows so the new
# dates render the same way as the existing ones.
$ws.Range("D89").Copy()
$ws.Range("D90:D94").PasteSpecial(-4122)

# Row 90
$ws.Range("B90").Value = "Selenium na notebooku, konzultace, mensi upravy"
$ws.Range("C90").Value = 4
$ws.Range("D90").Value = 41044

# Row 91
$ws.Range("B91").Value = "Unit testy, psani"
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 41045

# Row 92
$ws.Range("B92").Value = "Psani"
$ws.Range("C92").Value = 9
$ws.Range("D92").Value = 41046

# Row 93
$ws.Range("B93").Value = "Psani"
$ws.Range("C93").Value = 7
$ws.Range("D93").Value = 41048

# Row 94
$ws.Range("B94").Value = "Psani, bugy, dodelani funkcionality, tvorba pilotnich testu"
$ws.Range("C94").Value = 11
$ws.Range("D94").Value = 41049

# --- Update the view state: scroll so row 72 is at the top and select E77,
#     matching where the author ended up after adding the rows above. ---
$ws.Range("E77").Select()
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
